$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 35 (shifts nothing below it since it's the last row) so
# that it inherits the same row/column formatting as the existing table rows.
$ws.Rows("35").Insert() | Out-Null

$ws.Range("A35").Value = 567
$ws.Range("B35").Value = "Permutation in String"
$ws.Range("C35").Value = "Medium"
$ws.Range("D35").Value = "String,Sliding Window,Counter"
$ws.Range("E35").Value = 45701

# Move the active selection the way it ends up after adding a new row
# (previously B36, now one row further down at B37).
$ws.Range("B37").Select() | Out-Null
